$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 6.25
$ws.Range("M2").Value = 1.07
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 2.7
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.57
$ws.Range("V2").Value = 1.54

# Row 3
$ws.Range("G3").Value = 1.75
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 1.37
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("V3").Value = 1.63

# Row 4
$ws.Range("I4").Value = 2.15
$ws.Range("K4").Value = 1.95
$ws.Range("M4").Value = 1.07
$ws.Range("O4").Value = 1.41
$ws.Range("P4").Value = 2.62
$ws.Range("V4").Value = 1.63

# Row 5
$ws.Range("G5").Value = 1.55

# Row 6
$ws.Range("G6").Value = 5.25
$ws.Range("H6").Value = 3.4
$ws.Range("J6").Value = 5.5
$ws.Range("K6").Value = 2.05
$ws.Range("L6").Value = 2.4
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2.3
$ws.Range("R6").Value = 1.6
$ws.Range("S6").Value = 1.5
$ws.Range("T6").Value = 2.5
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("X6").Value = 23
$ws.Range("Y6").Value = 17
$ws.Range("Z6").Value = 51
$ws.Range("AA6").Value = 41
$ws.Range("AC6").Value = 7.5
$ws.Range("AE6").Value = 19
$ws.Range("AF6").Value = 67
$ws.Range("AH6").Value = 5.5
$ws.Range("AJ6").Value = 9
$ws.Range("AM6").Value = 34
$ws.Range("AN6").Value = 6.5
$ws.Range("AO6").Value = 29
$ws.Range("AQ6").Value = 101
$ws.Range("AR6").Value = 151
$ws.Range("AT6").Value = 2.5
$ws.Range("AU6").Value = 9
$ws.Range("AV6").Value = 67
$ws.Range("AW6").Value = 3.6
$ws.Range("AY6").Value = 23
$ws.Range("BA6").Value = 51
$ws.Range("BB6").Value = 201

# Row 7
$ws.Range("O7").Value = 1.44
$ws.Range("P7").Value = 2.63
$ws.Range("Q7").Value = 2.4
$ws.Range("R7").Value = 1.53
$ws.Range("V7").Value = 1.57

# Row 8
$ws.Range("G8").Value = 2.57
$ws.Range("H8").Value = 3.55
$ws.Range("I8").Value = 2.4
$ws.Range("J8").Value = 3.05
$ws.Range("L8").Value = 2.87
$ws.Range("U8").Value = 1.5
$ws.Range("V8").Value = 2.27
$ws.Range("W8").Value = 11.25
$ws.Range("X8").Value = 15
$ws.Range("Y8").Value = 9.75
$ws.Range("AA8").Value = 19
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 14.5
$ws.Range("AD8").Value = 7.3
$ws.Range("AE8").Value = 11.75
$ws.Range("AH8").Value = 11.25
$ws.Range("AJ8").Value = 9.25
$ws.Range("AK8").Value = 26
$ws.Range("AL8").Value = 17.5
$ws.Range("AN8").Value = 4.75
$ws.Range("AO8").Value = 13
$ws.Range("AP8").Value = 18
$ws.Range("AR8").Value = 70
$ws.Range("AU8").Value = 6.4
$ws.Range("AW8").Value = 4.55
$ws.Range("AX8").Value = 12
$ws.Range("AZ8").Value = 45

# Row 9
$ws.Range("V9").Value = 1.6

# Row 11
$ws.Range("K11").Value = 1.91

# Row 12
$ws.Range("G12").Value = 1.6
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 6.5
$ws.Range("AC12").Value = 6.5
$ws.Range("AH12").Value = 13

# Row 13
$ws.Range("Q13").Value = 1.98
$ws.Range("R13").Value = 1.88

